$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price figures stored as literal text (e.g. "45.145.53",
# "1.00"). Some of the refreshed prices look like plain numbers (e.g. "1.00",
# "0.500") which Excel would otherwise silently re-interpret as numeric values
# and normalize (dropping trailing zeros). Force those specific cells to a
# text format first so the literal string is preserved, matching how the
# rest of the workbook already stores this column.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range('D2').Value = '45.224.44'
$ws.Range('E2').Value = '  -3.22%  '
$ws.Range('D3').Value = '2.388.33'
$ws.Range('E3').Value = '  +5.45%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '292.91'
$ws.Range('E5').Value = '  -2.79%  '
$ws.Range('D6').Value = '94.04'
$ws.Range('E6').Value = '  -6.74%  '
$ws.Range('E7').Value = '  -0.80%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '0.500'
$ws.Range('E9').Value = '  -2.91%  '
$ws.Range('D10').Value = '34.20'
$ws.Range('E10').Value = '  -4.12%  '
$ws.Range('E11').Value = '  -0.59%  '
$ws.Range('E12').Value = '  -2.79%  '
$ws.Range('D13').Value = '0.104'
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('D14').Value = '2.756.00'
$ws.Range('E14').Value = '  +5.40%  '
$ws.Range('D15').Value = '2.387.79'
$ws.Range('E15').Value = '  +5.27%  '
$ws.Range('D16').Value = '14.04'
$ws.Range('E16').Value = '  +3.32%  '
$ws.Range('D17').Value = '0.823'
$ws.Range('E17').Value = '  +3.42%  '
$ws.Range('D18').Value = '45.209.51'
$ws.Range('E18').Value = '  -3.29%  '
$ws.Range('D19').Value = '12.41'
$ws.Range('E19').Value = '  -3.71%  '
$ws.Range('D20').Value = '0.0₃0931'
$ws.Range('E20').Value = '  +0.24%  '
$ws.Range('D21').Value = '6.07'
$ws.Range('E21').Value = '  +2.87%  '
$ws.Range('D22').Value = '66.50'
$ws.Range('E22').Value = '  +1.69%  '
$ws.Range('D23').Value = '237.94'
$ws.Range('E23').Value = '  -4.63%  '
$ws.Range('D24').Value = '2.76'
$ws.Range('E24').Value = '  -3.02%  '
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('E26').Value = '  +0.40%  '
$ws.Range('E27').Value = '  -1.03%  '
$ws.Range('D28').Value = '37.28'
$ws.Range('E28').Value = '  -13.31%  '
$ws.Range('D29').Value = '9.52'
$ws.Range('E29').Value = '  -2.37%  '
$ws.Range('E30').Value = '  +18.43%  '
$ws.Range('D31').Value = '20.91'
$ws.Range('E31').Value = '  +5.18%  '
$ws.Range('B32').Value = 'WEMIXToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D32').Value = '2.70'
$ws.Range('E32').Value = '  -3.03%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').Value = '147.01'
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('D34').Value = '5.39'
$ws.Range('E34').Value = '  -0.99%  '
$ws.Range('D35').Value = '0.0756'
$ws.Range('E35').Value = '  -2.01%  '
$ws.Range('D36').Value = '1.96'
$ws.Range('E36').Value = '  +13.71%  '
$ws.Range('D37').Value = '0.112'
$ws.Range('E37').Value = '  -2.21%  '
$ws.Range('E38').Value = '  -1.32%  '
$ws.Range('D39').Value = '14.44'
$ws.Range('E39').Value = '  -11.46%  '
$ws.Range('D40').Value = '3.68'
$ws.Range('E40').Value = '  -5.12%  '
$ws.Range('E41').Value = '  -1.75%  '
$ws.Range('D42').Value = '1.967.72'
$ws.Range('E42').Value = '  +8.69%  '
$ws.Range('D43').Value = '3.14'
$ws.Range('E43').Value = '  -2.19%  '
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('D45').Value = '88.35'
$ws.Range('E45').Value = '  -2.78%  '
$ws.Range('E46').Value = '  -14.28%  '
$ws.Range('D47').Value = '8.37'
$ws.Range('E47').Value = '  +7.33%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '14.96'
$ws.Range('E48').Value = '  +17.01%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '98.99'
$ws.Range('E49').Value = '  +5.49%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.624.78'
$ws.Range('E50').Value = '  +5.33%  '
$ws.Range('E51').Value = '  -3.71%  '
